$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 5.041100000000001
$ws.Range("C4").Value = -14.23989999999999
$ws.Range("C5").Value = -14.6331
$ws.Range("B6").Value = 8.929600000000011
$ws.Range("C6").Value = -11.28570000000001
$ws.Range("B7").Value = 5.212599999999997
$ws.Range("B8").Value = 4.885399999999999
$ws.Range("C8").Value = -11.2326
$ws.Range("B16").Value = 8.643800000000004
$ws.Range("C16").Value = -11.8591
$ws.Range("B20").Value = 5.622599999999995
$ws.Range("B21").Value = 5.184599999999998
$ws.Range("C22").Value = -11.10849999999999
